$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-11-07 Friday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2025-11-08 Saturday", 2)

# Update the division-problem table. Only rows 1, 5, 9, 13, 17 (of the 20
# table rows) carry text; the rest are blank spacer rows. Addressing cells
# directly by (row, column) avoids ambiguity from duplicate old values
# (e.g. "66÷4=16, 2" and "63÷7=9, 0" each occur twice in the source).
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("51÷9=5, 6", "86÷5=17, 1", "52÷3=17, 1", "53÷7=7, 4", "81÷6=13, 3")
    5  = @("51÷2=25, 1", "72÷4=18, 0", "45÷4=11, 1", "57÷8=7, 1", "33÷7=4, 5")
    9  = @("38÷6=6, 2", "82÷8=10, 2", "26÷6=4, 2", "41÷9=4, 5", "30÷4=7, 2")
    13 = @("49÷3=16, 1", "94÷3=31, 1", "51÷2=25, 1", "76÷9=8, 4", "93÷5=18, 3")
    17 = @("39÷9=4, 3", "48÷2=24, 0", "86÷6=14, 2", "47÷4=11, 3", "65÷7=9, 2")
}

foreach ($rowIndex in $newValues.Keys) {
    $rowVals = $newValues[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $rowVals[$col - 1]
    }
}
